$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Chapter 13 progress update: "Факт (кон дня)" (F12) and "Помидоры" (H12)
$ws.Range("F12").Value = 416
$ws.Range("H12").Value = 10
